$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume/1h change (E) columns for rows 2-50

# Row 2
$ws.Range("D2").Value = "28.206.18"
$ws.Range("E2").Value = "  +0.84%  "

# Row 3
$ws.Range("D3").Value = "1.880.14"
$ws.Range("E3").Value = "  +1.26%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.43%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5140"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.26%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3912"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.69%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08365"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.56%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.122"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.42%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.232"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.80%  "

# Row 13
$ws.Range("D13").Value = "1.890.10"
$ws.Range("E13").Value = "  +1.71%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.29%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.261"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.13%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.010"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001102"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.73%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.98%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06681"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.20%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.007"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.041"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.67%  "

# Row 23
$ws.Range("D23").Value = "28.238.00"
$ws.Range("E23").Value = "  +0.83%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.275"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.58%  "

# Row 26
$ws.Range("D26").Value = "2.095.20"
$ws.Range("E26").Value = "  +1.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.490"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.98%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.53%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.76%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1060"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.52%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.038"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.48%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.862"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.94%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.610"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.34%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.671"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02447"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.84%  "

# Row 37
$ws.Range("E37").Value = "  +0.71%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2189"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.202"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.22%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6508"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.28%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.07%  "

# Row 42
$ws.Range("E42").Value = "  -1.53%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6136"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.60%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.25%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.286"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.43%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.675"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.50%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.020"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.231"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.36%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.40%  "

# Row 51: coin replaced (Aave -> Cronos)
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06923"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.18%  "
